$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.046971140552575
$ws.Cells.Item(2, 4).Value = 1.052256501597609
$ws.Cells.Item(2, 5).Value = 1.054418954131038
$ws.Cells.Item(2, 6).Value = 1.064986505882286
$ws.Cells.Item(2, 9).Value = 1.043176314603966
$ws.Cells.Item(2, 10).Value = 1.052022694742469
$ws.Cells.Item(2, 11).Value = 1.055005736653696
$ws.Cells.Item(2, 12).Value = 1.057162217365531
$ws.Cells.Item(2, 13).Value = 1.067700963367275
$ws.Cells.Item(2, 14).Value = 1.021140948893315

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04811049665252
$ws.Cells.Item(3, 4).Value = 1.053121253672735
$ws.Cells.Item(3, 5).Value = 1.055403823078129
$ws.Cells.Item(3, 6).Value = 1.065991035406303
$ws.Cells.Item(3, 9).Value = 1.043438696995859
$ws.Cells.Item(3, 10).Value = 1.052809400490087
$ws.Cells.Item(3, 11).Value = 1.055683164831751
$ws.Cells.Item(3, 12).Value = 1.057959886995395
$ws.Cells.Item(3, 13).Value = 1.068520330886146
$ws.Cells.Item(3, 14).Value = 1.021406544510382

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.048847791319519
$ws.Cells.Item(4, 4).Value = 1.053680580245494
$ws.Cells.Item(4, 5).Value = 1.056041426064152
$ws.Cells.Item(4, 6).Value = 1.066641200836714
$ws.Cells.Item(4, 9).Value = 1.043606764293574
$ws.Cells.Item(4, 10).Value = 1.053317953688032
$ws.Cells.Item(4, 11).Value = 1.056120639761001
$ws.Cells.Item(4, 12).Value = 1.058475742665028
$ws.Cells.Item(4, 13).Value = 1.069050069263405
$ws.Cells.Item(4, 14).Value = 1.021578121104192

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.049157763821329
$ws.Cells.Item(5, 4).Value = 1.053915666745964
$ws.Cells.Item(5, 5).Value = 1.056309552471008
$ws.Cells.Item(5, 6).Value = 1.066914570288216
$ws.Cells.Item(5, 9).Value = 1.043677009830303
$ws.Cells.Item(5, 10).Value = 1.053531630260013
$ws.Cells.Item(5, 11).Value = 1.056304346647832
$ws.Cells.Item(5, 12).Value = 1.058692538911734
$ws.Cells.Item(5, 13).Value = 1.06927266404678
$ws.Cells.Item(5, 14).Value = 1.021650184440315

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.049209810397503
$ws.Cells.Item(6, 4).Value = 1.053955135604858
$ws.Cells.Item(6, 5).Value = 1.056354576650522
$ws.Cells.Item(6, 6).Value = 1.066960472555208
$ws.Cells.Item(6, 9).Value = 1.043688780319715
$ws.Cells.Item(6, 10).Value = 1.053567500508189
$ws.Cells.Item(6, 11).Value = 1.056335179665696
$ws.Cells.Item(6, 12).Value = 1.0587289358793
$ws.Cells.Item(6, 13).Value = 1.069310032397319
$ws.Cells.Item(6, 14).Value = 1.021662280233527

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.048851933132142
$ws.Cells.Item(7, 4).Value = 1.053683721696954
$ws.Cells.Item(7, 5).Value = 1.056045008475598
$ws.Cells.Item(7, 6).Value = 1.066644853456207
$ws.Cells.Item(7, 9).Value = 1.043607704528953
$ws.Cells.Item(7, 10).Value = 1.053320809313151
$ws.Cells.Item(7, 11).Value = 1.056123095277241
$ws.Cells.Item(7, 12).Value = 1.058478639780178
$ws.Cells.Item(7, 13).Value = 1.06905304400736
$ws.Cells.Item(7, 14).Value = 1.021579084283186

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.047356180258872
$ws.Cells.Item(8, 4).Value = 1.052548795180077
$ws.Cells.Item(8, 5).Value = 1.054751727420257
$ws.Cells.Item(8, 6).Value = 1.065325956484914
$ws.Cells.Item(8, 9).Value = 1.043265342263073
$ws.Cells.Item(8, 10).Value = 1.052288668775607
$ws.Cells.Item(8, 11).Value = 1.055234855980049
$ws.Cells.Item(8, 12).Value = 1.057431853540679
$ws.Cells.Item(8, 13).Value = 1.067977964950478
$ws.Cells.Item(8, 14).Value = 1.021230766429631

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.044720865492524
$ws.Cells.Item(9, 4).Value = 1.050547193540757
$ws.Cells.Item(9, 5).Value = 1.052475310932635
$ws.Cells.Item(9, 6).Value = 1.063003185146699
$ws.Cells.Item(9, 9).Value = 1.042648954595798
$ws.Cells.Item(9, 10).Value = 1.050466088748033
$ws.Cells.Item(9, 11).Value = 1.053663035493306
$ws.Cells.Item(9, 12).Value = 1.05558506446089
$ws.Cells.Item(9, 13).Value = 1.066080120166655
$ws.Cells.Item(9, 14).Value = 1.020614832341806

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.042964209649567
$ws.Cells.Item(10, 4).Value = 1.049211654093255
$ws.Cells.Item(10, 5).Value = 1.050959392722885
$ws.Cells.Item(10, 6).Value = 1.06145555043368
$ws.Cells.Item(10, 9).Value = 1.042229229412555
$ws.Cells.Item(10, 10).Value = 1.049248463377292
$ws.Cells.Item(10, 11).Value = 1.052610703714653
$ws.Cells.Item(10, 12).Value = 1.054352380997456
$ws.Cells.Item(10, 13).Value = 1.064812598515304
$ws.Cells.Item(10, 14).Value = 1.020202764325934

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.04220359999684
$ws.Cells.Item(11, 4).Value = 1.04863308310572
$ws.Cells.Item(11, 5).Value = 1.050303385019501
$ws.Cells.Item(11, 6).Value = 1.060785617158178
$ws.Cells.Item(11, 9).Value = 1.042045397833089
$ws.Cells.Item(11, 10).Value = 1.048720605864914
$ws.Cells.Item(11, 11).Value = 1.052153976960758
$ws.Cells.Item(11, 12).Value = 1.053818261058483
$ws.Cells.Item(11, 13).Value = 1.064263204464282
$ws.Cells.Item(11, 14).Value = 1.020023991758526

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.041921079774664
$ws.Cells.Item(12, 4).Value = 1.048418135055939
$ws.Cells.Item(12, 5).Value = 1.050059773773202
$ws.Cells.Item(12, 6).Value = 1.060536804540306
$ws.Cells.Item(12, 9).Value = 1.041976800958749
$ws.Cells.Item(12, 10).Value = 1.048524442981721
$ws.Cells.Item(12, 11).Value = 1.051984168952382
$ws.Cells.Item(12, 12).Value = 1.053619810891072
$ws.Cells.Item(12, 13).Value = 1.064059052517512
$ws.Cells.Item(12, 14).Value = 1.019957535904344

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.041981681153729
$ws.Cells.Item(13, 4).Value = 1.04846424399747
$ws.Cells.Item(13, 5).Value = 1.050112026534203
$ws.Cells.Item(13, 6).Value = 1.060590174315053
$ws.Cells.Item(13, 9).Value = 1.04199152941362
$ws.Cells.Item(13, 10).Value = 1.048566524803869
$ws.Cells.Item(13, 11).Value = 1.052020600559148
$ws.Cells.Item(13, 12).Value = 1.053662381582965
$ws.Cells.Item(13, 13).Value = 1.064102847539315
$ws.Cells.Item(13, 14).Value = 1.01997179325338

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.042180246710076
$ws.Cells.Item(14, 4).Value = 1.048615316269305
$ws.Cells.Item(14, 5).Value = 1.050283246830355
$ws.Cells.Item(14, 6).Value = 1.060765049609027
$ws.Cells.Item(14, 9).Value = 1.042039733992074
$ws.Cells.Item(14, 10).Value = 1.04870439288075
$ws.Cells.Item(14, 11).Value = 1.052139943829113
$ws.Cells.Item(14, 12).Value = 1.053801858213245
$ws.Cells.Item(14, 13).Value = 1.064246330884839
$ws.Cells.Item(14, 14).Value = 1.020018499552959

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.042302590011288
$ws.Cells.Item(15, 4).Value = 1.048708391448217
$ws.Cells.Item(15, 5).Value = 1.050388749150336
$ws.Cells.Item(15, 6).Value = 1.060872800080364
$ws.Cells.Item(15, 9).Value = 1.042069392864199
$ws.Cells.Item(15, 10).Value = 1.048789325607467
$ws.Cells.Item(15, 11).Value = 1.05221345403784
$ws.Cells.Item(15, 12).Value = 1.053887787175821
$ws.Cells.Item(15, 13).Value = 1.064334724770444
$ws.Cells.Item(15, 14).Value = 1.020047269985915

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.043014689434207
$ws.Cells.Item(16, 4).Value = 1.049250046186879
$ws.Cells.Item(16, 5).Value = 1.051002938104662
$ws.Cells.Item(16, 6).Value = 1.061500015968072
$ws.Cells.Item(16, 9).Value = 1.04224138572388
$ws.Cells.Item(16, 10).Value = 1.049283482502412
$ws.Cells.Item(16, 11).Value = 1.05264099282218
$ws.Cells.Item(16, 12).Value = 1.054387821181934
$ws.Cells.Item(16, 13).Value = 1.064849048419367
$ws.Cells.Item(16, 14).Value = 1.020214621614635

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.043461378945632
$ws.Cells.Item(17, 4).Value = 1.049589738703813
$ws.Cells.Item(17, 5).Value = 1.05138830825922
$ws.Cells.Item(17, 6).Value = 1.061893506187854
$ws.Cells.Item(17, 9).Value = 1.042348713357047
$ws.Cells.Item(17, 10).Value = 1.049593288483562
$ws.Cells.Item(17, 11).Value = 1.052908892718522
$ws.Cells.Item(17, 12).Value = 1.054701382835
$ws.Cells.Item(17, 13).Value = 1.065171523067421
$ws.Cells.Item(17, 14).Value = 1.020319504650606

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.04372192872195
$ws.Cells.Item(18, 4).Value = 1.049787849198914
$ws.Cells.Item(18, 5).Value = 1.051613126258258
$ws.Cells.Item(18, 6).Value = 1.062123041942459
$ws.Cells.Item(18, 9).Value = 1.042411114304922
$ws.Cells.Item(18, 10).Value = 1.049773933492791
$ws.Cells.Item(18, 11).Value = 1.053065051912149
$ws.Cells.Item(18, 12).Value = 1.054884243332351
$ws.Cells.Item(18, 13).Value = 1.065359564127141
$ws.Cells.Item(18, 14).Value = 1.020380647923133

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.043810770060495
$ws.Cells.Item(19, 4).Value = 1.0498553952772
$ws.Cells.Item(19, 5).Value = 1.051689789852953
$ws.Cells.Item(19, 6).Value = 1.06220131102287
$ws.Cells.Item(19, 9).Value = 1.042432357246138
$ws.Cells.Item(19, 10).Value = 1.04983551867623
$ws.Cells.Item(19, 11).Value = 1.053118280815018
$ws.Cells.Item(19, 12).Value = 1.054946588154287
$ws.Cells.Item(19, 13).Value = 1.065423672290842
$ws.Cells.Item(19, 14).Value = 1.020401490563579

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.043413453050053
$ws.Cells.Item(20, 4).Value = 1.049553295641827
$ws.Cells.Item(20, 5).Value = 1.05134695774768
$ws.Cells.Item(20, 6).Value = 1.06185128638745
$ws.Cells.Item(20, 9).Value = 1.04233721895384
$ws.Cells.Item(20, 10).Value = 1.049560055398221
$ws.Cells.Item(20, 11).Value = 1.052880160175135
$ws.Cells.Item(20, 12).Value = 1.054667744228478
$ws.Cells.Item(20, 13).Value = 1.065136930066532
$ws.Cells.Item(20, 14).Value = 1.020308255130472

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.042121773998777
$ws.Cells.Item(21, 4).Value = 1.048570830382476
$ws.Cells.Item(21, 5).Value = 1.050232825086684
$ws.Cells.Item(21, 6).Value = 1.060713552353899
$ws.Cells.Item(21, 9).Value = 1.042025547602987
$ws.Cells.Item(21, 10).Value = 1.048663796740818
$ws.Cells.Item(21, 11).Value = 1.052104804610341
$ws.Cells.Item(21, 12).Value = 1.053760787320711
$ws.Cells.Item(21, 13).Value = 1.064204080894842
$ws.Cells.Item(21, 14).Value = 1.020004747140937

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.041309667145687
$ws.Cells.Item(22, 4).Value = 1.047952878163097
$ws.Cells.Item(22, 5).Value = 1.049532668258113
$ws.Cells.Item(22, 6).Value = 1.059998389845872
$ws.Cells.Item(22, 9).Value = 1.041827772838969
$ws.Cells.Item(22, 10).Value = 1.048099744551327
$ws.Cells.Item(22, 11).Value = 1.051616385764342
$ws.Cells.Item(22, 12).Value = 1.053190233609111
$ws.Cells.Item(22, 13).Value = 1.063617084033906
$ws.Cells.Item(22, 14).Value = 1.01981362012703

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.041740178514073
$ws.Cells.Item(23, 4).Value = 1.048280488796139
$ws.Cells.Item(23, 5).Value = 1.049903802201792
$ws.Cells.Item(23, 6).Value = 1.060377494349165
$ws.Cells.Item(23, 9).Value = 1.041932788979137
$ws.Cells.Item(23, 10).Value = 1.048398810461765
$ws.Cells.Item(23, 11).Value = 1.051875393261496
$ws.Cells.Item(23, 12).Value = 1.053492724740764
$ws.Cells.Item(23, 13).Value = 1.063928307471952
$ws.Cells.Item(23, 14).Value = 1.01991496856928

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.043435108688642
$ws.Cells.Item(24, 4).Value = 1.049569762774316
$ws.Cells.Item(24, 5).Value = 1.051365642146023
$ws.Cells.Item(24, 6).Value = 1.061870363637305
$ws.Cells.Item(24, 9).Value = 1.042342413402228
$ws.Cells.Item(24, 10).Value = 1.049575072182881
$ws.Cells.Item(24, 11).Value = 1.052893143490605
$ws.Cells.Item(24, 12).Value = 1.054682944173891
$ws.Cells.Item(24, 13).Value = 1.065152561318034
$ws.Cells.Item(24, 14).Value = 1.020313338407008

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.04540211544123
$ws.Cells.Item(25, 4).Value = 1.051064857295627
$ws.Cells.Item(25, 5).Value = 1.053063521162587
$ws.Cells.Item(25, 6).Value = 1.063603522888727
$ws.Cells.Item(25, 9).Value = 1.042809857281421
$ws.Cells.Item(25, 10).Value = 1.05093772225521
$ws.Cells.Item(25, 11).Value = 1.054070173768106
$ws.Cells.Item(25, 12).Value = 1.056062766482025
$ws.Cells.Item(25, 13).Value = 1.066571162927034
$ws.Cells.Item(25, 14).Value = 1.020774320996714

